$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 66942.39
$ws.Range("I100").Value = 42064.2
$ws.Range("J100").Value = 191333.33
$ws.Range("K100").Value = 42064.2
$ws.Range("L100").Value = 191333.33
$ws.Range("M100").Value = -41523.2
$ws.Range("N100").Value = -192415.33
$ws.Range("H137").Value = 8195.303
$ws.Range("I137").Value = 8863.621
$ws.Range("K137").Value = 26590.863
$ws.Range("M137").Value = -24040.863

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2142.152
$ws.Range("I32").Value = 2148.6135
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 2148.6135
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1861.6135
$ws.Range("N32").Value = -2574
$ws.Range("H44").Value = 73999
$ws.Range("J44").Value = 73999
$ws.Range("L44").Value = 73999
$ws.Range("N44").Value = -74975
$ws.Range("H74").Value = 4231.143
$ws.Range("I74").Value = 2574.05
$ws.Range("J74").Value = 8373.875
$ws.Range("K74").Value = 2574.05
$ws.Range("L74").Value = 8373.875
$ws.Range("M74").Value = -1700.05
$ws.Range("N74").Value = -10121.875
$ws.Range("H77").Value = 4231.143
$ws.Range("I77").Value = 2574.05
$ws.Range("J77").Value = 8373.875
$ws.Range("K77").Value = 12870.25
$ws.Range("L77").Value = 41869.375
$ws.Range("M77").Value = -8502.25
$ws.Range("N77").Value = -50605.375
$ws.Range("H122").Value = 442489.72
$ws.Range("I122").Value = 3101.5
$ws.Range("K122").Value = 9304.5
$ws.Range("M122").Value = -6854.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H99").Value = 16089.76
$ws.Range("I99").Value = 16089.76
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 16089.76
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -14591.76
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 2348.4473
$ws.Range("I134").Value = 1674.6061
$ws.Range("K134").Value = 5023.8183
$ws.Range("M134").Value = -2488.8183

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3249.3333
$ws.Range("I31").Value = 1156.1428
$ws.Range("K31").Value = 1156.1428
$ws.Range("M31").Value = -861.1428000000001
$ws.Range("H34").Value = 3249.3333
$ws.Range("I34").Value = 1156.1428
$ws.Range("K34").Value = 1156.1428
$ws.Range("M34").Value = -954.1428000000001
$ws.Range("H138").Value = 22666.334
$ws.Range("J138").Value = 22666.334
$ws.Range("L138").Value = 22666.334
$ws.Range("N138").Value = -32946.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 565.3333
$ws.Range("I47").Value = 198.33333
$ws.Range("K47").Value = 594.99999
$ws.Range("M47").Value = -163.99999
$ws.Range("H80").Value = 600000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 600000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13148
$ws.Range("I43").Value = 15900
$ws.Range("J43").Value = 10396
$ws.Range("K43").Value = 15900
$ws.Range("L43").Value = 10396
$ws.Range("M43").Value = -15749
$ws.Range("N43").Value = -10698
$ws.Range("H46").Value = 7000
$ws.Range("I46").Value = 7000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -6844
$ws.Range("N46").ClearContents()
$ws.Range("H52").Value = 30666.666
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30666.666
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 30666.666
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -31184.666
$ws.Range("H102").Value = 5576.968
$ws.Range("J102").Value = 2187.5
$ws.Range("L102").Value = 2187.5
$ws.Range("N102").Value = -5431.5
$ws.Range("H113").Value = 2100.6
$ws.Range("I113").Value = 1625.75
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1625.75
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 544.25
$ws.Range("N113").Value = -8340
$ws.Range("H131").Value = 53326
$ws.Range("J131").Value = 53326
$ws.Range("L131").Value = 53326
$ws.Range("N131").Value = -63406
$ws.Range("H132").Value = 2326.7837
$ws.Range("I132").Value = 2483.7812
$ws.Range("J132").Value = 1322
$ws.Range("K132").Value = 7451.3436
$ws.Range("L132").Value = 3966
$ws.Range("M132").Value = -4921.3436
$ws.Range("N132").Value = -9026

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2930.0557
$ws.Range("I22").Value = 3826
$ws.Range("J22").Value = 2034.1111
$ws.Range("K22").Value = 3826
$ws.Range("L22").Value = 2034.1111
$ws.Range("M22").Value = -3531
$ws.Range("N22").Value = -2624.1111
$ws.Range("H27").Value = 2930.0557
$ws.Range("I27").Value = 3826
$ws.Range("J27").Value = 2034.1111
$ws.Range("K27").Value = 3826
$ws.Range("L27").Value = 2034.1111
$ws.Range("M27").Value = -3719
$ws.Range("N27").Value = -2248.1111
$ws.Range("H40").Value = 14663
$ws.Range("I40").Value = 16346
$ws.Range("K40").Value = 16346
$ws.Range("M40").Value = -16210
$ws.Range("H61").Value = 4448.758
$ws.Range("I61").Value = 2932.5715
$ws.Range("K61").Value = 2932.5715
$ws.Range("M61").Value = -2730.5715
$ws.Range("H93").Value = 6163.346
$ws.Range("I93").Value = 6881.1577
$ws.Range("J93").Value = 4215
$ws.Range("K93").Value = 6881.1577
$ws.Range("L93").Value = 4215
$ws.Range("M93").Value = -5633.1577
$ws.Range("N93").Value = -6711
$ws.Range("H113").Value = 4448.758
$ws.Range("I113").Value = 2932.5715
$ws.Range("K113").Value = 2932.5715
$ws.Range("M113").Value = -762.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16612.5
$ws.Range("I41").Value = 20000
$ws.Range("J41").Value = 13225
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 13225
$ws.Range("M41").Value = -19610
$ws.Range("N41").Value = -14005
$ws.Range("H107").Value = 34253.473
$ws.Range("I107").Value = 3546.3076
$ws.Range("J107").Value = 100785.664
$ws.Range("K107").Value = 10638.9228
$ws.Range("L107").Value = 302356.992
$ws.Range("M107").Value = -8718.9228
$ws.Range("N107").Value = -306196.992
$ws.Range("H131").Value = 41650
$ws.Range("I131").Value = 41650
$ws.Range("K131").Value = 41650
$ws.Range("M131").Value = -36610
$ws.Range("H132").Value = 18017.648
$ws.Range("I132").Value = 22681.809
$ws.Range("J132").Value = 6993.273
$ws.Range("K132").Value = 68045.427
$ws.Range("L132").Value = 20979.819
$ws.Range("M132").Value = -65515.427
$ws.Range("N132").Value = -26039.819
$ws.Range("H141").Value = 79398.6
$ws.Range("J141").Value = 79398.6
$ws.Range("L141").Value = 79398.6
$ws.Range("N141").Value = -89758.6
